# Fix alignment issue: a new row was inserted at the top of the data
# (row 2), shifting the existing data rows down by one and dropping the
# final (now duplicate/overflow) row so the sheet keeps the same extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2; this pushes the existing data rows
# (2-11) down to (3-12).
$ws.Rows("2:2").Insert()

# Inserting a row carries over formatting from the row above (the bold /
# centered / bordered header style) which is not wanted here - the data
# rows have no explicit style. Strip that back off.
$ws.Range("A2:D2").ClearFormats()

# Populate the newly inserted row with the new record.
$ws.Range("A2").Value = "Move Robot32 to location (2, 9) and remove the toolkit."
$ws.Range("B2").Value = 66.63938
$ws.Range("C2").Value = 9970

# Column D stores the cost as text (e.g. "0.02688"), not a number, so
# force a text number format before assigning to avoid Excel silently
# converting the numeric-looking string into a real number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.02688"
$ws.Range("D2").ClearFormats()

# The insert shifted the old last data row (previously row 11, the
# "Move to location (9, 5) and remove the toolkit." entry) down to row
# 12, which now falls outside the table - remove it entirely so the
# sheet's dimensions/content match the target.
$ws.Rows("12:12").Delete()
